$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8741305470466614
$ws.Range("B1").Value = 0.6899652481079102
$ws.Range("C1").Value = 3.94738507270813
$ws.Range("D1").Value = 2.984235048294067
$ws.Range("E1").Value = 0.8164118528366089
